$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Quantity for row 7 (Hitec HS-5645MG) changed from 1 to 2.
# This also recalculates the Total in E7 (formula C7*D7) and the
# grand total sum in E26 (SUM(E4:E22)) automatically.
$ws.Range("C7").Value = 2

# New line item added in row 23: "Arduino Due"
$ws.Range("B23").Value = "Arduino Due"

# Update the view so the active cell / selection matches the new state,
# and scroll the window so the newly added row is visible.
$ws.Range("B22").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
